# TC05_Canine_Filter_Diagnosis-OsteoSarcoma.xlsx
#
# The FilesTab Neo4j query (stored in B4 of the "startup" sheet) is
# rewritten to drop the `File Type` and `Breed` columns from its RETURN
# clause (those fields aren't available off the `file` node for this
# query shape). Everything else on the sheet is left untouched; Excel's
# own shared-string table compaction/re-indexing happens automatically
# when the cell content changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFilesTabQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.disease_term IN ['Osteosarcoma']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFilesTabQuery

# Match the author's final cursor position/selection on the sheet.
$ws.Range("B4").Select()
